$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update labels to clarify these are US indicators
$ws.Range("B2").Value = "US Real GDP"
$ws.Range("B3").Value = "US Unemployment"

# Re-apply merges on column A groupings (unmerge then re-merge)
$ws.Range("A2:A3").UnMerge()
$ws.Range("A4:A5").UnMerge()
$ws.Range("A6:A8").UnMerge()
$ws.Range("A10:A11").UnMerge()
$ws.Range("A12:A14").UnMerge()
$ws.Range("A15:A17").UnMerge()

$ws.Range("A12:A14").Merge()
$ws.Range("A4:A5").Merge()
$ws.Range("A15:A17").Merge()
$ws.Range("A2:A3").Merge()
$ws.Range("A10:A11").Merge()
$ws.Range("A6:A8").Merge()
